# The DiSCoVER "top drugs (cerebellar stem cell control)" results table is
# re-emitted by the notebook pipeline as a new trailing slide (slide 33),
# identical in content to the existing copies of that slide already in the
# deck (e.g. slide 32). Duplicate that slide; PowerPoint appends the copy
# immediately after the source, which lands it at the end of the deck since
# slide 32 is currently last.
$p = $ppt.ActivePresentation
$src = $p.Slides.Item($p.Slides.Count)
$src.Duplicate() | Out-Null
